# "Add cycle time numbers medium"
# Adds a row of "cycle" header-like numbers (21,24,33,41,48,56) across E2:J2
# on Sheet2, matching the styling already used in columns A/B, widens the
# E:J columns, and updates several existing rows (22,25,34,42,49,57) whose
# draw counts/dates rolled forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- New block E2:J2 -------------------------------------------------
# Values first.
$ws.Range("E2").Value = 21
$ws.Range("F2").Value = 24
$ws.Range("G2").Value = 33
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 56

# Styling: reuse the same cell formats already present in column B so the
# resulting style indices line up (B2 uses the "red" style, B3 uses the
# plain-border style).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null

$ws.Range("B3").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Re-assert the values (PasteSpecial of formats only shouldn't disturb
# them, but make sure nothing got clobbered).
$ws.Range("E2").Value = 21
$ws.Range("F2").Value = 24
$ws.Range("G2").Value = 33
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 56

# Widen the new columns E:J to match column A's display width (3 chars).
# Excel's ColumnWidth property and the stored OOXML <col width> differ by a
# constant ~0.8333 padding factor, so compensate to land exactly on "3".
$ws.Columns("E:J").ColumnWidth = 2.1666666666666665

# --- Row 22: draw count + date rolled forward -------------------------
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 45291

# --- Row 25: now "hit" (1) instead of "miss" (0), restyled yellow -----
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 45291

# --- Row 34: same pattern as row 25 -----------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 45291

# --- Row 42: draw count + date rolled forward -------------------------
$ws.Range("B42").Value = 4
$ws.Range("C42").Value = 45291

# --- Row 49: now "hit" (1) instead of "miss" (0), restyled yellow -----
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A49").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B49").Value = 1
$ws.Range("C49").Value = 45291

# --- Row 57: draw count + date rolled forward -------------------------
$ws.Range("B57").Value = 3
$ws.Range("C57").Value = 45291

# --- Selection the author left on the sheet ---------------------------
$ws.Range("J6").Select() | Out-Null
